$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "birthPlace" (F) and fastestDriver tuple (G) columns entirely.
$ws.Range("F1:G5").EntireColumn.Delete()

# Row 2 & 3: birthDate (E) was blank -> now "nan"
$ws.Range("E2").Value = "nan"
$ws.Range("E3").Value = "nan"

# Row 4: replace with the corrected record (George Murray, bishop of Rochester)
$ws.Range("A4").Value = "http://dbpedia.org/resource/George_Murray_(bishop_of_Rochester)"
$ws.Range("B4").Value = "http://dbpedia.org/resource/Chester_Square"
$ws.Range("C4").Value = "http://dbpedia.org/resource/Lord_George_Murray_(bishop)"
$ws.Range("D4").Value = "http://dbpedia.org/resource/1860"
$ws.Range("E4").Value = "nan"

# Row 5: replace with the corrected record (Geoffrey, archbishop of York)
$ws.Range("A5").Value = "http://dbpedia.org/resource/Geoffrey_(archbishop_of_York)"
$ws.Range("B5").Value = "http://dbpedia.org/resource/Normandy"
$ws.Range("C5").Value = "http://dbpedia.org/resource/Henry_II_of_England"
$ws.Range("D5").Value = "http://dbpedia.org/resource/12-12-12"
$ws.Range("E5").Value = "nan"

# Row 6: Ferdinand III of Castile
$ws.Range("A6").Value = "http://dbpedia.org/resource/Ferdinand_III_of_Castile"
$ws.Range("B6").Value = "http://dbpedia.org/resource/Crown_of_Castile"
$ws.Range("C6").Value = "http://dbpedia.org/resource/Alfonso_IX_of_León"
$ws.Range("D6").Value = "http://dbpedia.org/resource/1252"
$ws.Range("E6").Value = "nan"

# Row 7: Erik Benzelius the younger
$ws.Range("A7").Value = "http://dbpedia.org/resource/Erik_Benzelius_the_younger"
$ws.Range("B7").Value = "http://dbpedia.org/resource/Linköping"
$ws.Range("C7").Value = "http://dbpedia.org/resource/Erik_Benzelius_the_Elder"
$ws.Range("D7").Value = "http://dbpedia.org/resource/1743"
$ws.Range("E7").Value = "nan"

# Row 8: Edward the Confessor
$ws.Range("A8").Value = "http://dbpedia.org/resource/Edward_the_Confessor"
$ws.Range("B8").Value = "http://dbpedia.org/resource/London"
$ws.Range("C8").Value = "http://dbpedia.org/resource/Æthelred_the_Unready"
$ws.Range("D8").Value = "http://dbpedia.org/resource/1066"
$ws.Range("E8").Value = "nan"

# Row 9: Edward William Grinfield
$ws.Range("A9").Value = "http://dbpedia.org/resource/Edward_William_Grinfield"
$ws.Range("B9").Value = "http://dbpedia.org/resource/Brighton"
$ws.Range("C9").Value = "http://dbpedia.org/resource/Thomas_Grinfield"
$ws.Range("D9").Value = "http://dbpedia.org/resource/1864"
$ws.Range("E9").Value = "nan"

# Row 10: Edward Francis Wilson
$ws.Range("A10").Value = "http://dbpedia.org/resource/Edward_Francis_Wilson"
$ws.Range("B10").Value = "http://dbpedia.org/resource/Saltspring_Island"
$ws.Range("C10").Value = "http://dbpedia.org/resource/Daniel_Wilson_(bishop)"
$ws.Range("D10").Value = "http://dbpedia.org/resource/1915"
$ws.Range("E10").Value = "nan"

# Row 11: Donald Foster Hudson
$ws.Range("A11").Value = "http://dbpedia.org/resource/Donald_Foster_Hudson"
$ws.Range("B11").Value = "http://dbpedia.org/resource/England"
$ws.Range("C11").Value = "http://dbpedia.org/resource/Father"
$ws.Range("D11").Value = "http://dbpedia.org/resource/2003"
$ws.Range("E11").Value = "nan"
